# Rename the repeated-measures time-point labels in row 1 (B1:G1) on every
# sheet of the workbook: baseline/retest/follow-up -> 1baseline/2retest/3followup
# for both the EXP and CON groups.

$wb = $excel.ActiveWorkbook

$renames = @{
    "EXP_baseline"   = "EXP_1baseline"
    "EXP_retest"     = "EXP_2retest"
    "EXP_follow-up"  = "EXP_3followup"
    "CON_baseline"   = "CON_1baseline"
    "CON_retest"     = "CON_2retest"
    "CON_follow-up"  = "CON_3followup"
}

foreach ($ws in $wb.Worksheets) {
    foreach ($col in 2..7) {
        $cell = $ws.Cells.Item(1, $col)
        $current = $cell.Value2
        if ($renames.ContainsKey($current)) {
            $cell.Value = $renames[$current]
        }
    }
}
